$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All "Price" (column D) values are numeric-looking text (e.g. "232.04",
# "0.000009180") that must be stored as literal text, matching the source
# data (inline strings), not auto-coerced into floating point numbers by
# Excel's smart-entry parsing (which would also lose exact formatting like
# trailing zeros, e.g. "2.500" -> 2.5). Force text storage by switching
# NumberFormat to Text before each write, then snap the style back to the
# default "Normal" cell style (clears the explicit s= index) so the cells
# end up identical in formatting to the originals.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'

$ws.Range('D2').Value = '28.658.60'
$ws.Range('E2').Value = '  -1.86%  '
$ws.Range('D3').Value = '1.802.00'
$ws.Range('E3').Value = '  -1.31%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').Value = '232.04'
$ws.Range('E5').Value = '  -0.82%  '
$ws.Range('D6').Value = '0.5914'
$ws.Range('E6').Value = '  -1.40%  '
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('E8').Value = '  -0.59%  '
$ws.Range('D9').Value = '0.06807'
$ws.Range('E9').Value = '  -3.50%  '
$ws.Range('D10').Value = '23.28'
$ws.Range('E10').Value = '  -0.60%  '
$ws.Range('D11').Value = '0.07507'
$ws.Range('E11').Value = '  -1.75%  '
$ws.Range('D12').Value = '1.798.24'
$ws.Range('E12').Value = '  -1.60%  '
$ws.Range('D13').Value = '4.763'
$ws.Range('E13').Value = '  -0.29%  '
$ws.Range('E14').Value = '  -0.46%  '
$ws.Range('D15').Value = '2.046.76'
$ws.Range('D16').Value = '0.000009180'
$ws.Range('E16').Value = '  -7.33%  '
$ws.Range('D17').Value = '75.58'
$ws.Range('E17').Value = '  -4.27%  '
$ws.Range('D18').Value = '28.645.07'
$ws.Range('E18').Value = '  -1.76%  '
$ws.Range('D19').Value = '5.460'
$ws.Range('E19').Value = '  -6.24%  '
$ws.Range('D20').Value = '1.004'
$ws.Range('D21').Value = '210.48'
$ws.Range('E21').Value = '  -6.76%  '
$ws.Range('D22').Value = '11.49'
$ws.Range('D23').Value = '6.826'
$ws.Range('E23').Value = '  -2.32%  '
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('D25').Value = '153.85'
$ws.Range('E25').Value = '  -0.99%  '
$ws.Range('D26').Value = '7.845'
$ws.Range('E26').Value = '  -1.88%  '
$ws.Range('D27').Value = '0.1268'
$ws.Range('E27').Value = '  -2.13%  '
$ws.Range('D28').Value = '16.46'
$ws.Range('E28').Value = '  -0.44%  '
$ws.Range('D29').Value = '1.415'
$ws.Range('E29').Value = '  -4.28%  '
$ws.Range('D30').Value = '0.06150'
$ws.Range('E30').Value = '  -0.41%  '
$ws.Range('E31').Value = '  -1.11%  '
$ws.Range('D32').Value = '3.771'
$ws.Range('E32').Value = '  -1.29%  '
$ws.Range('D33').Value = '3.739'
$ws.Range('E33').Value = '  -1.34%  '
$ws.Range('E34').Value = '  -0.81%  '
$ws.Range('E35').Value = '  -5.61%  '
$ws.Range('D36').Value = '0.6416'
$ws.Range('E36').Value = '  +0.38%  '
$ws.Range('D37').Value = '2.500'
$ws.Range('E37').Value = '  -1.36%  '
$ws.Range('D38').Value = '2.715'
$ws.Range('E38').Value = '  -0.62%  '
$ws.Range('D39').Value = '6.539'
$ws.Range('E39').Value = '  +0.13%  '
$ws.Range('E40').Value = '  -2.88%  '
$ws.Range('D41').Value = '1.149.81'
$ws.Range('E41').Value = '  -5.49%  '
$ws.Range('D42').Value = '0.8837'
$ws.Range('E42').Value = '  -1.97%  '
$ws.Range('D43').Value = '1.006'
$ws.Range('E43').Value = '  +0.48%  '
$ws.Range('D44').Value = '100.08'
$ws.Range('E44').Value = '  -0.34%  '
$ws.Range('D45').Value = '1.952.26'
$ws.Range('D46').Value = '60.37'
$ws.Range('E46').Value = '  -3.51%  '
$ws.Range('E47').Value = '  -3.17%  '
$ws.Range('D48').Value = '1.583'
$ws.Range('E48').Value = '  +0.43%  '
$ws.Range('D49').Value = '8.371'
$ws.Range('E49').Value = '  -1.25%  '
$ws.Range('E50').Value = '  -0.73%  '
$ws.Range('E51').Value = '  -1.69%  '

$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
